$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.346.21"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.905.69"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("E5").Value = "  +9.57%  "
$ws.Range("D6").Value = "'246.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'41.59"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.63%  "
$ws.Range("E9").Value = "  +4.51%  "
$ws.Range("E10").Value = "  +11.90%  "
$ws.Range("D11").Value = "'0.0727"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").Value = "'0.0995"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "2.181.47"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "'12.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").Value = "1.907.52"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "'4.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "35.298.48"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'72.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "0.0₃0823"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "'241.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").Value = "'12.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").Value = "'4.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "'2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E26").Value = "  +8.59%  "
$ws.Range("D27").Value = "'168.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("E28").Value = "  +2.47%  "
$ws.Range("E29").Value = "  +4.36%  "
$ws.Range("D30").Value = "'18.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("D33").Value = "'0.963"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "'0.0574"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'4.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +8.79%  "
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "'0.0664"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.44%  "
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("D43").Value = "'16.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.88%  "
$ws.Range("D44").Value = "'90.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "1.349.03"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "'2.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.36%  "
$ws.Range("D47").Value = "'12.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Value = "'46.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("E51").Value = "  -2.22%  "
